$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=133; A="2025-08-05 11:34"; B="Hacer Yaren Ünsal"; C="Hacer Yaren Ünsal"; D=24; E=26; F=20; G=26; H=25; I=23; J=0.25; K=0.27; L=0.23; M=0.25; N="%24.96"; O="%27.13"; P="%22.96"; Q="%24.96" },
    @{ Row=134; A="2025-08-05 11:44"; B="Hacer Yaren Ünsal"; C="Hacer Yaren Ünsal"; D=24; E=22; F=16; G=21; H=23; I=18.5; J=0.18; K=0.29; L=0.2; M=0.32; N="%18.47"; O="%29.45"; P="%20.07"; Q="%32.01" },
    @{ Row=135; A="2025-08-05 11:51"; B="Hacer Yaren Ünsal"; C="Hacer Yaren Ünsal"; D=24; E=33; F=29; G=29; H=28.5; I=29; J=0.36; K=0.24; L=0.25; M=0.16; N="%35.87"; O="%23.5"; P="%24.54"; Q="%16.08" },
    @{ Row=136; A="2025-08-07 13:58"; B="yaren"; C="yaren"; D=24; E=21; F=23; G=18; H=22.5; I=20.5; J=0.2; K=0.27; L=0.23; M=0.3; N="%20.02"; O="%26.86"; P="%22.69"; Q="%30.44" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
}
